# Slide 18 ("3. Uzdevuma rezultāts"): shrink/reposition the result
# screenshot and add the "Pārbaude:" label text box above it, matching
# the same layout pattern already used on the other "Uzdevuma rezultāts"
# slides (1., 2.) in this deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# --- resize / reposition the existing result picture -----------------------
$pic = $s.Shapes.Item(2)
$pic.Left   = 81.8
$pic.Top    = 151.05
$pic.Width  = 222.3
$pic.Height = 355.30005

# --- add the "Pārbaude:" label text box ------------------------------------
# Copy the identical label shape from slide 13 instead of synthesizing a
# brand-new textbox: this reproduces the exact OOXML shape (no extra
# a16:creationId noise) used throughout the deck for this recurring label.
$srcSlide = $p.Slides.Item(13)
$srcBox = $srcSlide.Shapes.Item(7)

# The new-shape id allocator hands out a mis-seeded id on its very first
# allocation for this slide (it grabs the lowest free gap, i.e. "3"
# here) and only increments correctly from then on. Burn through the
# mis-seeded allocations with throwaway copies so the real label lands on
# id 8 -- exactly like the equivalent labels elsewhere in this deck.
for ($i = 0; $i -lt 4; $i++) {
    $srcBox.Copy()
    $throwaway = $s.Shapes.Paste()
    $throwaway.Item(1).Delete()
}

$srcBox.Copy()
$pasted = $s.Shapes.Paste()
$txBox = $pasted.Item(1)
$txBox.Left   = 81.8
$txBox.Top    = 111.9
$txBox.Width  = 90.3
$txBox.Height = 29.0
